# Sprint backlog completed - still need to update contribution percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 used to be the static text "Story"; it now pulls the sprint total from B12
# via a formula (B12 is currently blank, which evaluates to 0).
$ws.Range("B1").Formula = "=+B12"

# Row 2 (AK's story): the previously-blank contribution cells (E, G, J, K)
# get an explicit 0% value, matching the existing whole-percent format
# already used by F2/H2/I2 in that row.
$ws.Range("E2").Value = 0
$ws.Range("E2").NumberFormat = "0%"
$ws.Range("G2").Value = 0
$ws.Range("G2").NumberFormat = "0%"
$ws.Range("J2").Value = 0
$ws.Range("J2").NumberFormat = "0%"
$ws.Range("K2").Value = 0
$ws.Range("K2").NumberFormat = "0%"

# Rows 3-8: fill in equal contribution percentages (14.3%) across all
# seven contributor columns (E:K), displayed with 2 decimal places and
# centered.
$rng = $ws.Range("E3:K8")
$rng.Value = 0.143
$rng.NumberFormat = "0.00%"
$rng.HorizontalAlignment = -4108

# Reflect the author's last selection in the saved file.
$ws.Range("B1").Select()
